$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-number cells in column D to retain string type (not auto-converted to numeric)
$textCells = 'D5','D6','D14','D19','D21','D25','D26','D27','D32','D33','D34','D35','D39','D40','D41','D42','D43','D45','D46','D47','D48','D50'
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply cell value updates per the diff
$ws.Range('D2').Value = '63.171.86'
$ws.Range('E2').Value = '  +1.07%  '
$ws.Range('D3').Value = '2.449.65'
$ws.Range('E3').Value = '  +0.66%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').Value = '571.75'
$ws.Range('E5').Value = '  +0.99%  '
$ws.Range('D6').Value = '146.34'
$ws.Range('E6').Value = '  +0.80%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  +0.93%  '
$ws.Range('D9').Value = '2.446.11'
$ws.Range('E9').Value = '  +0.56%  '
$ws.Range('E10').Value = '  +0.63%  '
$ws.Range('E11').Value = '  +1.36%  '
$ws.Range('E12').Value = '  -0.35%  '
$ws.Range('E13').Value = '  +0.33%  '
$ws.Range('D14').Value = '26.99'
$ws.Range('E14').Value = '  +0.79%  '
$ws.Range('E15').Value = '  +0.26%  '
$ws.Range('D16').Value = '2.895.13'
$ws.Range('E16').Value = '  +1.05%  '
$ws.Range('D17').Value = '63.050.17'
$ws.Range('E17').Value = '  +1.05%  '
$ws.Range('D18').Value = '2.457.05'
$ws.Range('E18').Value = '  +0.17%  '
$ws.Range('D19').Value = '11.28'
$ws.Range('E19').Value = '  +0.45%  '
$ws.Range('E20').Value = '  +5.20%  '
$ws.Range('D21').Value = '328.39'
$ws.Range('E21').Value = '  +1.55%  '
$ws.Range('E22').Value = '  +0.91%  '
$ws.Range('E23').Value = '  +13.83%  '
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('D25').Value = '65.69'
$ws.Range('E25').Value = '  -2.27%  '
$ws.Range('D26').Value = '615.39'
$ws.Range('E26').Value = '  +5.35%  '
$ws.Range('D27').Value = '8.92'
$ws.Range('E27').Value = '  +4.29%  '
$ws.Range('E28').Value = '  +2.77%  '
$ws.Range('D29').Value = '2.563.35'
$ws.Range('E29').Value = '  +0.31%  '
$ws.Range('E30').Value = '  +4.84%  '
$ws.Range('E31').Value = '  +0.12%  '
$ws.Range('D32').Value = '8.25'
$ws.Range('E32').Value = '  -2.02%  '
$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').Value = '1.91'
$ws.Range('E33').Value = '  +1.68%  '
$ws.Range('B34').Value = 'Kaspa'
$ws.Range('C34').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D34').Value = '0.141'
$ws.Range('E34').Value = '  -3.17%  '
$ws.Range('D35').Value = '5.18'
$ws.Range('E35').Value = '  +7.28%  '
$ws.Range('E36').Value = '  +1.79%  '
$ws.Range('E37').Value = '  +0.05%  '
$ws.Range('E38').Value = '  -0.38%  '
$ws.Range('D39').Value = '5.42'
$ws.Range('E39').Value = '  +1.82%  '
$ws.Range('D40').Value = '18.81'
$ws.Range('E40').Value = '  +0.47%  '
$ws.Range('D41').Value = '147.01'
$ws.Range('E41').Value = '  -0.56%  '
$ws.Range('D42').Value = '1.78'
$ws.Range('E42').Value = '  -1.92%  '
$ws.Range('D43').Value = '2.59'
$ws.Range('E43').Value = '  +7.11%  '
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('D45').Value = '41.84'
$ws.Range('E45').Value = '  +0.58%  '
$ws.Range('D46').Value = '148.50'
$ws.Range('E46').Value = '  +0.25%  '
$ws.Range('D47').Value = '3.76'
$ws.Range('E47').Value = '  +2.82%  '
$ws.Range('D48').Value = '21.16'
$ws.Range('E48').Value = '  +3.53%  '
$ws.Range('E49').Value = '  -0.10%  '
$ws.Range('D50').Value = '0.601'
$ws.Range('E50').Value = '  +0.03%  '
$ws.Range('E51').Value = '  +0.65%  '
